$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26, shifting existing rows 26-47 down to 27-48
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new weekly record
$ws.Cells.Item(26, 1).Value = 10
$ws.Cells.Item(26, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(26, 3).Value = 'La Araucanía'
$ws.Cells.Item(26, 4).Value = 44452
$ws.Cells.Item(26, 5).Value = 9
$ws.Cells.Item(26, 6).Value = 100112035
$ws.Cells.Item(26, 7).Value = 'Bruselas (repollito)'
$ws.Cells.Item(26, 8).Value = 'Sin especificar'
$ws.Cells.Item(26, 9).Value = 'Primera'
$ws.Cells.Item(26, 10).Value = 80
$ws.Cells.Item(26, 11).Value = 25000
$ws.Cells.Item(26, 12).Value = 25000
$ws.Cells.Item(26, 13).Value = 25000
$ws.Cells.Item(26, 14).Value = '$/malla 10 kilos'
$ws.Cells.Item(26, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(26, 16).Value = 2500
$ws.Cells.Item(26, 17).Value = 10
$ws.Cells.Item(26, 18).Value = 'Hortaliza'
